$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 5.8
$ws.Range("G2").Value = 6.4
$ws.Range("H2").Value = 1.77
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.52
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 2.8
$ws.Range("O2").Value = 1.48
$ws.Range("P2").Value = 1.61
$ws.Range("Q2").Value = 2.44
$ws.Range("R2").Value = 1.22
$ws.Range("S2").Value = 4.8
$ws.Range("T2").Value = 2.26
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 2.18
$ws.Range("W2").Value = 1.17
$ws.Range("Y2").Value = 6.6
$ws.Range("AB2").Value = 16
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 11.5
$ws.Range("AG2").Value = 25
$ws.Range("AH2").Value = 28
$ws.Range("AK2").Value = 960
$ws.Range("AL2").Value = 960
$ws.Range("AN2").Value = 200
$ws.Range("AO2").Value = 22
# Row 3
$ws.Range("F3").Value = 1.01
$ws.Range("H3").Value = 1.01
$ws.Range("N3").Value = 1.36
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 1.36
$ws.Range("Q3").Value = 1.14
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 1.03
# Row 4
$ws.Range("F4").Value = 6.4
$ws.Range("H4").Value = 1.59
$ws.Range("I4").Value = 1.62
$ws.Range("J4").Value = 4.2
$ws.Range("K4").Value = 4.5
$ws.Range("L4").Value = 1.45
$ws.Range("N4").Value = 3.2
$ws.Range("O4").Value = 1.39
$ws.Range("P4").Value = 1.78
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.28
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 2.18
$ws.Range("U4").Value = 1.74
$ws.Range("V4").Value = 2.3
$ws.Range("X4").Value = 16.5
$ws.Range("Y4").Value = 6.8
$ws.Range("Z4").Value = 8.4
$ws.Range("AA4").Value = 15
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 19
$ws.Range("AG4").Value = 27
$ws.Range("AH4").Value = 30
$ws.Range("AI4").Value = 46
$ws.Range("AO4").Value = 12
# Row 5
$ws.Range("F5").Value = 1.01
$ws.Range("H5").Value = 1.01
$ws.Range("J5").Value = 1.03
$ws.Range("O5").Value = 1.3
$ws.Range("Q5").Value = 1.3
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 1.03
$ws.Range("U5").Value = 1.03
# Row 6
$ws.Range("I6").Value = 2.4
$ws.Range("V6").Value = 1.71
# Row 7
$ws.Range("F7").Value = 1.68
$ws.Range("G7").Value = 1.71
$ws.Range("I7").Value = 7.4
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 3.85
$ws.Range("N7").Value = 3.2
$ws.Range("O7").Value = 1.41
$ws.Range("P7").Value = 1.75
$ws.Range("S7").Value = 4
$ws.Range("W7").Value = 2.4
$ws.Range("X7").Value = 11.5
$ws.Range("Y7").Value = 19
$ws.Range("AA7").Value = 220
$ws.Range("AI7").Value = 130
$ws.Range("AN7").Value = 13
# Row 8
$ws.Range("F8").Value = 2.14
$ws.Range("G8").Value = 2.28
$ws.Range("H8").Value = 4.1
$ws.Range("I8").Value = 4.6
$ws.Range("M8").Value = 1.15
$ws.Range("N8").Value = 2.4
$ws.Range("O8").Value = 1.67
$ws.Range("P8").Value = 1.45
$ws.Range("Q8").Value = 3
$ws.Range("S8").Value = 6.6
$ws.Range("T8").Value = 2.42
$ws.Range("U8").Value = 1.58
$ws.Range("V8").Value = 1.27
$ws.Range("W8").Value = 1.78
$ws.Range("Y8").Value = 970
$ws.Range("Z8").Value = 34
$ws.Range("AB8").Value = 6.4
$ws.Range("AE8").Value = 95
$ws.Range("AH8").Value = 32
$ws.Range("AJ8").Value = 44
$ws.Range("AK8").Value = 46
# Row 9
$ws.Range("F9").Value = 1.97
$ws.Range("G9").Value = 2.06
$ws.Range("H9").Value = 4.7
$ws.Range("I9").Value = 5.3
$ws.Range("K9").Value = 3.45
$ws.Range("T9").Value = 2.2
$ws.Range("U9").Value = 1.65
$ws.Range("V9").Value = 1.23
$ws.Range("W9").Value = 1.94
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 36
$ws.Range("AA9").Value = 160
$ws.Range("AB9").Value = 6.8
$ws.Range("AD9").Value = 23
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 150
$ws.Range("AJ9").Value = 25
$ws.Range("AK9").Value = 29
$ws.Range("AM9").Value = 290
$ws.Range("AN9").Value = 25
# Row 10
$ws.Range("G10").Value = 2.58
$ws.Range("H10").Value = 3.3
$ws.Range("I10").Value = 3.5
$ws.Range("N10").Value = 2.84
$ws.Range("Q10").Value = 2.46
$ws.Range("T10").Value = 2.02
$ws.Range("W10").Value = 1.63
$ws.Range("Y10").Value = 11.5
$ws.Range("Z10").Value = 27
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 9.199999999999999
$ws.Range("AC10").Value = 7.6
$ws.Range("AD10").Value = 16.5
$ws.Range("AE10").Value = 160
$ws.Range("AF10").Value = 16.5
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 24
$ws.Range("AJ10").Value = 40
$ws.Range("AK10").Value = 36
$ws.Range("AL10").Value = 1000
# Row 11
$ws.Range("F11").Value = 1.79
$ws.Range("G11").Value = 1.82
$ws.Range("H11").Value = 5.1
$ws.Range("I11").Value = 5.7
$ws.Range("O11").Value = 1.4
$ws.Range("V11").Value = 1.21
$ws.Range("W11").Value = 2.22
$ws.Range("Y11").Value = 16
$ws.Range("AL11").Value = 46
# Row 12
$ws.Range("I12").Value = 4.9
$ws.Range("N12").Value = 2.76
$ws.Range("P12").Value = 1.59
$ws.Range("Q12").Value = 2.52
$ws.Range("R12").Value = 1.21
$ws.Range("S12").Value = 5.1
$ws.Range("U12").Value = 1.77
$ws.Range("X12").Value = 9.6
$ws.Range("AE12").Value = 1000
# Row 13
$ws.Range("F13").Value = 1.49
$ws.Range("G13").Value = 1.51
$ws.Range("J13").Value = 5
$ws.Range("O13").Value = 1.23
$ws.Range("P13").Value = 2.32
$ws.Range("Q13").Value = 1.69
$ws.Range("R13").Value = 1.54
$ws.Range("S13").Value = 2.68
$ws.Range("T13").Value = 1.83
$ws.Range("U13").Value = 2.08
$ws.Range("W13").Value = 2.92
$ws.Range("X13").Value = 30
$ws.Range("Y13").Value = 27
$ws.Range("Z13").Value = 210
$ws.Range("AA13").Value = 960
$ws.Range("AC13").Value = 11
$ws.Range("AD13").Value = 28
$ws.Range("AE13").Value = 960
$ws.Range("AF13").Value = 9.800000000000001
$ws.Range("AI13").Value = 85
$ws.Range("AM13").Value = 960
$ws.Range("AN13").Value = 6.4
$ws.Range("AO13").Value = 100
